$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# The sheet holds DDAIF's yearly financials laid out with "Period Ending"
# header rows (7, 38, 80) followed by line-item rows, one column per
# fiscal year (most recent first). This update adds the newly reported
# FY2018 (period ending 2018-12-31) column in front of the existing
# FY2017..FY2011 columns (which shift one column to the right), and
# restates the figures for the previously-reported years.
# -----------------------------------------------------------------------

# Step 1: insert a new blank column at D; old D:K shifts right to E:L
$ws.Columns("D:D").Insert()

# Step 2: give the new column D the same number formats/fonts as column E
# (bold date style on the "Period Ending" rows, right-aligned "#,##0" on
# the data rows) by copying the formatting across.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Step 3: match the new column's width to the other year columns
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# Target values for D:K (8 cols) per row, after the new column is inserted
$data = @{
  7 = @(43465, 43100, 42735, 42369, 42004, 41639, 41274, 40908)
  8 = @(187778500, 184179100, 171957300, 167700500, 145715100, 132374600, 128240100, 125057700)
  9 = @(150528400, 145439100, 135952700, 132362300, 113986300, 104102700, 103853600, 190211200)
  10 = @(37250100, 38740100, 36004700, 35338200, 31728800, 28271900, 24386500, -65153500)
  11 = @("", "", "", "", "", "", "", "")
  12 = @(7309800, 6608500, 5855700, 5298000, 5036600, 4696700, 9329300, 4899500)
  13 = @(0, 0, 0, 0, 0, 0, 0, 0)
  14 = @(149200, 28000, 309700, 340000, 273800, 158200, 151500, 0)
  15 = @(233400, 204200, 164900, 170500, 182900, 68400, 70700, "NA")
  16 = @("", "", "", "", "", "", "", "")
  17 = @(176277000, 169540500, 158366600, 153407400, 135181800, 123606300, 119169900, 115140200)
  18 = @(11501500, 14638600, 13590700, 14293000, 10533200, 8768400, 9070200, 9917500)
  19 = @("", "", "", "", "", "", "", "")
  20 = @(1126500, 1448500, 875200, 352300, 1290300, 3201000, 699000, 272300)
  21 = @(19719200, 22470800, 20626900, 20700700, 17445900, 16882000, 14343300, 14389100)
  22 = @(740500, 416300, 357900, 346700, 409500, 593500, 663100, 272300)
  23 = @(11887500, 15670800, 14107900, 14298600, 11414000, 11375900, 9106100, 9917500)
  24 = @(3380600, 3758700, 4252300, 4525000, 3234700, 1592100, 1442900, 2840600)
  25 = @(0, 0, 0, 0, 0, 0, 0, 0)
  26 = @(8506900, 11912200, 9855600, 9773700, 8179300, 9783800, 7663200, 7076900)
  27 = @(8133300, 11531800, 9566100, 9451600, 7811300, 7676700, 7212200, 6652000)
  28 = @(0, 0, 0, 0, 0, 0, 0, 0)
  29 = @(0, 0, 0, 0, 0, 0, 0, 0)
  30 = @(0, 0, 0, 0, 0, 0, 0, 0)
  31 = @(0, 0, 0, 0, 0, 0, 0, 0)
  32 = @(-1126500, -1448500, -875200, -352300, -1290300, -3201000, -699000, -272300)
  33 = @(8133300, 11531800, 9566100, 9451600, 7811300, 7676700, 7212200, 6652000)
  34 = @(0, 0, 0, 0, 0, 0, 0, 0)
  35 = @(8133300, 11531800, 9566100, 9451600, 7811300, 7676700, 7212200, 6652000)
  38 = @(43465, 43100, 42735, 42369, 42004, 41639, 41274, 40908)
  39 = @("", "", "", "", "", "", "", "")
  40 = @("", "", "", "", "", "", "", "")
  41 = @(17786900, 13544700, 12320600, 11148100, 10846300, 12401400, 12337400, 11240400)
  42 = @(12679600, 16877000, 13275400, 10405300, 7897700, 7979600, 6876700, 1565900)
  43 = @(69741800, 74933200, 58207700, 53070100, 42759000, 35291100, 66289400, 33346800)
  44 = @(33086400, 28819400, 28480600, 26658500, 23409200, 19465400, 19881700, 20049800)
  45 = @(3153900, 8887300, 2217100, 1769400, 1643700, 3896700, 6187800, 5538000)
  46 = @(136449000, 119936000, 114501000, 103051000, 86555900, 79034100, 75687200, 71740900)
  47 = @(66787600, 65249300, 56088300, 52913000, 46620900, 40086500, 43010400, 35936200)
  48 = @(90234900, 168422000, 82267700, 70981600, 63091700, 56031100, 104697000, 49289500)
  49 = @(16606600, 30821100, 13573800, 11297300, 10509700, 10533200, 19937800, 9694500)
  50 = @(0, 0, 0, 0, 0, 0, 0, 0)
  51 = @(0, 0, 0, 0, 0, 0, 0, 0)
  52 = @(5896100, 6632100, 6199000, 5414700, 5990300, 3390700, 8030100, 7217800)
  53 = @(0, 0, 0, 0, 0, 0, 0, 0)
  54 = @(315974000, 286495000, 272630000, 243658000, 212769000, 189076000, 182954000, 173879000)
  55 = @("", "", "", "", "", "", "", "")
  56 = @("", "", "", "", "", "", "", "")
  57 = @(15915400, 13969900, 12978100, 11834800, 11419600, 10194400, 9909400, 11168800)
  58 = @(63100700, 54692500, 53056700, 46350500, 40717000, 37016700, 36987500, 31341900)
  59 = @(30885000, 40134700, 28725200, 28298800, 23007500, 19107500, 20883600, 21878600)
  60 = @(109901000, 98313300, 94759900, 86484100, 75144200, 66318600, 65878800, 64389300)
  61 = @(99477900, 87939300, 78985900, 67129800, 56547200, 50204600, 48627000, 41630300)
  62 = @(32483900, 27134200, 32537700, 28756600, 31054400, 23899500, 24320300, 19337300)
  63 = @(0, 0, 0, 0, 0, 0, 0, 0)
  64 = @(0, 0, 0, 0, 0, 0, 0, 0)
  65 = @(0, 0, 0, 0, 0, 0, 0, 0)
  66 = @(243418000, 214834000, 207611000, 183563000, 163777000, 141189000, 140425000, 127368000)
  67 = @("", "", "", "", "", "", "", "")
  68 = @(0, 0, 0, 0, 0, 0, 0, 0)
  69 = @(0, 0, 0, 0, 0, 0, 0, 0)
  70 = @(0, 0, 0, 0, 0, 0, 0, 0)
  71 = @(0, 0, 0, 0, 0, 0, 0, 0)
  72 = @(55527300, 53354000, 45770500, 41503500, 31962100, 30998300, 38815200, 41687900)
  73 = @(0, 0, 0, 0, 0, 0, 0, 0)
  74 = @(0, 0, 0, 0, 0, 0, 0, 0)
  75 = @(0, 0, 0, 0, 0, 0, 0, 0)
  76 = @(72555700, 71660400, 65019300, 60094900, 48991700, 47886500, 42529000, 46511000)
  77 = @(0, 0, 0, 0, 0, 0, 0, 0)
  80 = @(43465, 43100, 42735, 42369, 42004, 41639, 41274, 40908)
  81 = @(8133300, 11531800, 9566100, 9451600, 7811300, 7676700, 7212200, 6652000)
  82 = @("", "", "", "", "", "", "", "")
  83 = @(7074100, 6368400, 6146300, 6040800, 5608800, 4900900, 4563100, 4196400)
  84 = @(0, 0, 0, 0, 0, 0, 0, 0)
  85 = @(0, 0, 0, 0, 0, 0, 0, 0)
  86 = @(0, 0, 0, 0, 0, 0, 0, 0)
  87 = @(0, 0, 0, 0, 0, 0, 0, 0)
  88 = @(0, 0, 0, 0, 0, 0, 0, 0)
  89 = @(384800, -1853500, 4163700, 249100, -1429400, 3685700, -1234200, -817000)
  90 = @("", "", "", "", "", "", "", "")
  91 = @(-8453100, -7566700, -6607400, -5694100, -5434900, -5581900, -5415800, -6897300)
  92 = @(0, 0, 0, 0, 0, 0, 0, 0)
  93 = @(0, 0, 0, 0, 0, 0, 0, 0)
  94 = @(-11131300, -10679100, -16455100, -10908000, -3039500, -7662100, -9945300, -7673200)
  95 = @("", "", "", "", "", "", "", "")
  96 = @(-4381400, -3901200, -3901200, -2940700, -2700600, -2635600, -2632200, -2313600)
  97 = @(0, 0, 0, 0, 0, 0, 0, 0)
  98 = @(0, 0, 0, 0, 0, 0, 0, 0)
  99 = @(0, 0, 0, 0, 0, 0, 0, 0)
  100 = @(14839400, 14730600, 13474000, 10805900, 2551400, 4325300, 12909600, 6857400)
  101 = @(149200, -973900, -10100, 154800, 362400, -285000, -136900, 75100)
  102 = @(4242200, 1224100, 1172500, 301800, -1555100, 64000, 1593200, -1557600)
}
# Step 4: write the restated values (columns D..K, left to right) for
# every row that carries year-by-year figures.
foreach ($entry in $data.GetEnumerator()) {
    $r = [int]$entry.Key
    $vals = $entry.Value
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $col = 4 + $i   # D=4 ... K=11
        $ws.Cells.Item($r, $col).Value2 = $vals[$i]
    }
}
